$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32 (Item ID 5484)
$ws.Range("H32").Value = 2200.2
$ws.Range("J32").Value = 2000.6666
$ws.Range("L32").Value = 2000.6666
$ws.Range("N32").Value = -2652.6666
# Row 58 (Item ID 4606)
$ws.Range("H58").Value = 1209.9333
$ws.Range("I58").Value = 469.44446
$ws.Range("J58").Value = 2320.6667
$ws.Range("K58").Value = 1408.33338
$ws.Range("L58").Value = 6962.000100000001
$ws.Range("M58").Value = -1258.33338
$ws.Range("N58").Value = -7262.000100000001
# Row 64 (Item ID 5506)
$ws.Range("H64").Value = 3044.5
$ws.Range("I64").Value = 2712.5
$ws.Range("K64").Value = 2712.5
$ws.Range("M64").Value = -2464.5
# Row 67 (Item ID 5506)
$ws.Range("H67").Value = 3044.5
$ws.Range("I67").Value = 2712.5
$ws.Range("K67").Value = 2712.5
$ws.Range("M67").Value = -1854.5
# Row 100 (Item ID 19906)
$ws.Range("H100").Value = 2767.5
$ws.Range("I100").Value = 2767.5
$ws.Range("K100").Value = 2767.5
$ws.Range("M100").Value = -2226.5
# Row 116 (Item ID 27778)
$ws.Range("H116").Value = 16333.223
$ws.Range("I116").Value = 52500
$ws.Range("K116").Value = 52500
$ws.Range("M116").Value = -49058
# Row 137 (Item ID 44013)
$ws.Range("H137").Value = 2107.5386
$ws.Range("I137").Value = 1750
$ws.Range("J137").Value = 2414
$ws.Range("K137").Value = 5250
$ws.Range("L137").Value = 7242
$ws.Range("M137").Value = -2700
$ws.Range("N137").Value = -12342
# Row 141 (Item ID 44161)
$ws.Range("H141").Value = 3169.0715
$ws.Range("J141").Value = 5634.8335
$ws.Range("L141").Value = 16904.5005
$ws.Range("N141").Value = -27264.5005

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Item ID 44147)
$ws.Range("H32").Value = 2461.747
$ws.Range("I32").Value = 1660.0139
$ws.Range("K32").Value = 1660.0139
$ws.Range("M32").Value = -1373.0139
# Row 63 (Item ID 12528)
$ws.Range("H63").Value = 6919.6
$ws.Range("I63").Value = 6919.6
$ws.Range("K63").Value = 6919.6
$ws.Range("M63").Value = -6233.6
# Row 66 (Item ID 12528)
$ws.Range("H66").Value = 6919.6
$ws.Range("I66").Value = 6919.6
$ws.Range("K66").Value = 34598
$ws.Range("M66").Value = -31166
# Row 135 (Item ID 42016)
$ws.Range("H135").Value = 19429
$ws.Range("J135").Value = 19429
$ws.Range("L135").Value = 19429
$ws.Range("N135").Value = -29569

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Item ID 14149)
$ws.Range("H20").Value = 1881.2727
$ws.Range("I20").Value = 2026.7778
$ws.Range("K20").Value = 2026.7778
$ws.Range("M20").Value = -1779.7778
# Row 105 (Item ID 19947)
$ws.Range("H105").Value = 2519.1538
$ws.Range("I105").Value = 2479.0833
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2479.0833
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -732.0832999999998
$ws.Range("N105").Value = -6494

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Item ID 44023)
$ws.Range("H31").Value = 1740.409
$ws.Range("I31").Value = 1405.8125
$ws.Range("J31").Value = 2632.6667
$ws.Range("K31").Value = 1405.8125
$ws.Range("L31").Value = 2632.6667
$ws.Range("M31").Value = -1110.8125
$ws.Range("N31").Value = -3222.6667
# Row 34 (Item ID 44023)
$ws.Range("H34").Value = 1740.409
$ws.Range("I34").Value = 1405.8125
$ws.Range("J34").Value = 2632.6667
$ws.Range("K34").Value = 1405.8125
$ws.Range("L34").Value = 2632.6667
$ws.Range("M34").Value = -1203.8125
$ws.Range("N34").Value = -3036.6667
# Row 44 (Item ID 1850)
$ws.Range("H44").Value = 6640
$ws.Range("I44").Value = 3000
$ws.Range("J44").Value = 10280
$ws.Range("K44").Value = 3000
$ws.Range("L44").Value = 10280
$ws.Range("N44").Value = -11164
# Row 62 (Item ID 12580)
$ws.Range("H62").Value = 2710.8
$ws.Range("I62").Value = 2651
$ws.Range("J62").Value = 2950
$ws.Range("K62").Value = 2651
$ws.Range("L62").Value = 2950
$ws.Range("M62").Value = -2027
$ws.Range("N62").Value = -4198
# Row 65 (Item ID 12580)
$ws.Range("H65").Value = 2710.8
$ws.Range("I65").Value = 2651
$ws.Range("J65").Value = 2950
$ws.Range("K65").Value = 13255
$ws.Range("L65").Value = 14750
$ws.Range("M65").Value = -10135
$ws.Range("N65").Value = -20990
# Row 107 (Item ID 27689)
$ws.Range("H107").Value = 541.53845
$ws.Range("I107").Value = 453.33334
$ws.Range("K107").Value = 453.33334
$ws.Range("M107").Value = 1466.66666
# Row 132 (Item ID 44019)
$ws.Range("H132").Value = 2920.5
$ws.Range("I132").Value = 1714
$ws.Range("J132").Value = 3570.1538
$ws.Range("K132").Value = 5142
$ws.Range("L132").Value = 10710.4614
$ws.Range("M132").Value = -2612
$ws.Range("N132").Value = -15770.4614

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 39 (Item ID 4712)
$ws.Range("H39").Value = 2500
$ws.Range("J39").Value = 2500
$ws.Range("L39").Value = 7500
$ws.Range("N39").Value = -8088
# Row 131 (Item ID 36060)
$ws.Range("H131").Value = 1240.24
$ws.Range("I131").Value = 609.5
$ws.Range("J131").Value = 1280.5
$ws.Range("K131").Value = 1828.5
$ws.Range("L131").Value = 3841.5
$ws.Range("M131").Value = 3211.5
$ws.Range("N131").Value = -13921.5
# Row 138 (Item ID 44105)
$ws.Range("H138").Value = 3344.2666
$ws.Range("I138").Value = 2444.8572
$ws.Range("J138").Value = 4131.25
$ws.Range("K138").Value = 7334.571599999999
$ws.Range("L138").Value = 12393.75
$ws.Range("M138").Value = -2194.571599999999
$ws.Range("N138").Value = -22673.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Item ID 14146)
$ws.Range("H70").Value = 4643.625
$ws.Range("I70").Value = 4999.5
$ws.Range("K70").Value = 4999.5
$ws.Range("M70").Value = -4729.5
# Row 73 (Item ID 14146)
$ws.Range("H73").Value = 4643.625
$ws.Range("I73").Value = 4999.5
$ws.Range("K73").Value = 4999.5
$ws.Range("M73").Value = -4063.5
# Row 107 (Item ID 27802)
$ws.Range("H107").Value = 1034.3334
$ws.Range("I107").Value = 500
$ws.Range("K107").Value = 500
$ws.Range("M107").Value = 1420

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Item ID 36249)
$ws.Range("H7").Value = 3017.5
$ws.Range("J7").Value = 4161.5
$ws.Range("L7").Value = 4161.5
$ws.Range("N7").Value = -4385.5
# Row 126 (Item ID 36249)
$ws.Range("H126").Value = 3017.5
$ws.Range("J126").Value = 4161.5
$ws.Range("L126").Value = 12484.5
$ws.Range("N126").Value = -17424.5
# Row 136 (Item ID 44060)
$ws.Range("H136").Value = 4166.154
$ws.Range("I136").Value = 2961.6667
$ws.Range("K136").Value = 8885.000100000001
$ws.Range("M136").Value = -6335.000100000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 48 (Item ID 3140)
$ws.Range("H48").Value = 100059
$ws.Range("I48").Value = 100059
$ws.Range("K48").Value = 100059
$ws.Range("M48").Value = -99490
# Row 107 (Item ID 27746)
$ws.Range("H107").Value = 1334.3334
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1501.5
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 4504.5
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -8344.5
# Row 132 (Item ID 44029)
$ws.Range("H132").Value = 1255.5151
$ws.Range("I132").Value = 966.375
$ws.Range("K132").Value = 2899.125
$ws.Range("M132").Value = -369.125
